$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny float-precision correction to the existing 15:00 timestamp
$ws.Cells.Item(15, 1).Value = 45878.62517811343

# Add new row of sensor data for 16:00:18
$ws.Cells.Item(16, 1).Value = 45878.66688301708
$ws.Cells.Item(16, 2).Value = 2025
$ws.Cells.Item(16, 3).Value = 37
$ws.Cells.Item(16, 4).Value = 18.22
$ws.Cells.Item(16, 5).Value = 78.92
$ws.Cells.Item(16, 6).Value = 293.55
$ws.Cells.Item(16, 7).Value = 17.36
$ws.Cells.Item(16, 8).Value = "ESE"
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = "16:00:18"

# Match the date/time number format used by the preceding rows in column A
$ws.Cells.Item(16, 1).NumberFormat = $ws.Cells.Item(15, 1).NumberFormat
